$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(111, "customername", "Customer's Name", "internal", "pii"),
    @(111555, "customername", "Customer's Name", "internal", "pii"),
    @(2222, "customername", "Customer's Name", "internal", "pii"),
    @(22333332, "customername", "Customer's Name", "internal", "pii"),
    @(22, "customername", "Customer's Name", "internal", "pii"),
    @(224, "customername", "Customer's Name", "internal", "pii"),
    @(33, "customername", "Customer's Name", "internal", "pii"),
    @(10000, "customername", "Customer's Name", "internal", "pii")
)

$startRow = 12
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
}
